# Updated C3DC phs002371 queries
# - TreatmentTab query (cell B5) gets an extra filter condition
#   (AND trt.treatment_id IS NOT NULL) added to the WHERE clause.
# - The sheet's active selection moves from B5 to C5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$treatmentCell = $ws.Range("B5")
$query = $treatmentCell.Value()

$oldClause = "std.dbgap_accession = 'phs002371' AND srv.cause_of_death IN ('Not Reported')`nORDER BY"
$newClause = "std.dbgap_accession = 'phs002371' AND srv.cause_of_death IN ('Not Reported') AND trt.treatment_id IS NOT NULL`nORDER BY"

if ($query.Contains($oldClause)) {
    $query = $query.Replace($oldClause, $newClause)
}

$treatmentCell.Value = $query

# Move the active selection to C5 (was B5).
$ws.Range("C5").Select()
